# Refresh market-data-derived columns (currentAveragePrice*, LevePrice*, LeveProfit*)
# on the per-job Leve tables, pulled from this runs market snapshot.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 40: 'Stuck in the Moment' / 'Horn Glue'
$ws.Range("H40").Value = 1913.4872
$ws.Range("I40").Value = 1879.2963
$ws.Range("J40").Value = 1990.4166
$ws.Range("K40").Value = 1879.2963
$ws.Range("L40").Value = 1990.4166
$ws.Range("M40").Value = -1704.2963
$ws.Range("N40").Value = -2340.4166

# Row 62: 'The Mustache Suits Him' / 'Enchanted Mythrite Ink'
$ws.Range("H62").Value = 2320.4614
$ws.Range("I62").Value = 1361
$ws.Range("K62").Value = 1361
$ws.Range("M62").Value = -737

# Row 65: 'Forgery of Convenience (L)' / 'Enchanted Mythrite Ink'
$ws.Range("H65").Value = 2320.4614
$ws.Range("I65").Value = 1361
$ws.Range("K65").Value = 6805
$ws.Range("M65").Value = -3685

# Row 69: 'Steeling the Knife, Steeling the Mind' / 'Grade 1 Mind Dissolvent'
$ws.Range("H69").Value = 3685.75
$ws.Range("I69").Value = 2256.5
$ws.Range("J69").Value = 5115
$ws.Range("K69").Value = 6769.5
$ws.Range("L69").Value = 15345
$ws.Range("M69").Value = -5895.5
$ws.Range("N69").Value = -17093

# Row 72: 'Surgical Substitution (L)' / 'Grade 1 Mind Dissolvent'
$ws.Range("H72").Value = 3685.75
$ws.Range("I72").Value = 2256.5
$ws.Range("J72").Value = 5115
$ws.Range("K72").Value = 20308.5
$ws.Range("L72").Value = 46035
$ws.Range("M72").Value = -15940.5
$ws.Range("N72").Value = -54771

# Row 113: 'Amaro Kart' / 'Starch Glue'
$ws.Range("H113").Value = 5531.2
$ws.Range("I113").Value = 4538.25
$ws.Range("K113").Value = 4538.25
$ws.Range("M113").Value = -1284.25

# Row 121: 'Mindful Medicine' / 'Tincture of Mind'
$ws.Range("H121").Value = 1157.5883
$ws.Range("I121").Value = 566.6667
$ws.Range("K121").Value = 1700.0001
$ws.Range("M121").Value = 46.99990000000003

$ws = $wb.Worksheets.Item("ARM")
# Row 27: 'Get Me the Hard Stuff' / 'Ironclad Bronze Buckler'
$ws.Range("H27").Value = 84990
$ws.Range("J27").Value = 84990
$ws.Range("L27").Value = 84990
$ws.Range("N27").Value = -85358

# Row 43: 'They''ve Got Legs' / 'Steel Sabatons'
$ws.Range("H43").Value = 12094.25
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 12094.25
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 12094.25
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -12720.25

# Row 92: 'Mail It In' / 'High Steel Scale Mail of Fending'
$ws.Range("H92").Value = 29000
$ws.Range("J92").Value = 29000
$ws.Range("L92").Value = 29000
$ws.Range("N92").Value = -33992

# Row 97: 'Ore for Me' / 'High Steel Ingot'
$ws.Range("H97").Value = 522.7273
$ws.Range("I97").Value = 462.77777
$ws.Range("J97").Value = 792.5
$ws.Range("K97").Value = 462.77777
$ws.Range("L97").Value = 792.5
$ws.Range("M97").Value = 33.22223000000002
$ws.Range("N97").Value = -1784.5

# Row 110: 'Scheduled Maintenance' / 'Deepgold Ingot'
$ws.Range("H110").Value = 1474.4546
$ws.Range("I110").Value = 599.1905
$ws.Range("J110").Value = 3006.1667
$ws.Range("K110").Value = 599.1905
$ws.Range("L110").Value = 3006.1667
$ws.Range("M110").Value = 1445.8095
$ws.Range("N110").Value = -7096.1667

$ws = $wb.Worksheets.Item("BSM")
# Row 20: 'Smelt and Dealt' / 'Iron Ingot'
$ws.Range("H20").Value = 2821.3157
$ws.Range("I20").Value = 2555.9092
$ws.Range("J20").Value = 3186.25
$ws.Range("K20").Value = 2555.9092
$ws.Range("L20").Value = 3186.25
$ws.Range("M20").Value = -2308.9092
$ws.Range("N20").Value = -3680.25

$ws = $wb.Worksheets.Item("CRP")
# Row 31: 'Wall Not Found' / 'Walnut Lumber'
$ws.Range("H31").Value = 1967.75
$ws.Range("I31").Value = 1214.16
$ws.Range("J31").Value = 2721.34
$ws.Range("K31").Value = 1214.16
$ws.Range("L31").Value = 2721.34
$ws.Range("M31").Value = -919.1600000000001
$ws.Range("N31").Value = -3311.34

# Row 34: 'Armoires of the Rich and Famous' / 'Walnut Lumber'
$ws.Range("H34").Value = 1967.75
$ws.Range("I34").Value = 1214.16
$ws.Range("J34").Value = 2721.34
$ws.Range("K34").Value = 1214.16
$ws.Range("L34").Value = 2721.34
$ws.Range("M34").Value = -1012.16
$ws.Range("N34").Value = -3125.34

# Row 50: 'The Arsenal of Theocracy' / 'Cobalt Halberd'
$ws.Range("H50").Value = 9491
$ws.Range("J50").Value = 9491
$ws.Range("L50").Value = 9491
$ws.Range("N50").Value = -10741

# Row 51: 'Greenstone for Greenhorns' / 'Jade Crook'
$ws.Range("H51").Value = 23082
$ws.Range("J51").Value = 23082
$ws.Range("L51").Value = 23082
$ws.Range("N51").Value = -24554

# Row 61: 'Incant Now, Think Later' / 'Jade Crook'
$ws.Range("H61").Value = 23082
$ws.Range("J61").Value = 23082
$ws.Range("L61").Value = 23082
$ws.Range("N61").Value = -23778

# Row 62: 'Splinter in the Sewers' / 'Cedar Lumber'
$ws.Range("H62").Value = 3265.25
$ws.Range("I62").Value = 2362
$ws.Range("J62").Value = 4770.6665
$ws.Range("K62").Value = 2362
$ws.Range("L62").Value = 4770.6665
$ws.Range("M62").Value = -1738
$ws.Range("N62").Value = -6018.6665

# Row 65: 'The Lumber of Their Discontent (L)' / 'Cedar Lumber'
$ws.Range("H65").Value = 3265.25
$ws.Range("I65").Value = 2362
$ws.Range("J65").Value = 4770.6665
$ws.Range("K65").Value = 11810
$ws.Range("L65").Value = 23853.3325
$ws.Range("M65").Value = -8690
$ws.Range("N65").Value = -30093.3325

# Row 68: 'Do You Even String Bow' / 'Holy Cedar Composite Bow'
$ws.Range("H68").Value = 29530
$ws.Range("J68").Value = 29530
$ws.Range("L68").Value = 29530
$ws.Range("N68").Value = -31028

# Row 71: 'Win One Bow, Get Three Free (L)' / 'Holy Cedar Composite Bow'
$ws.Range("H71").Value = 29530
$ws.Range("J71").Value = 29530
$ws.Range("L71").Value = 88590
$ws.Range("N71").Value = -96078

# Row 74: 'License to Heal' / 'Dark Chestnut Rod'
$ws.Range("H74").Value = 16861.46
$ws.Range("J74").Value = 16861.46
$ws.Range("L74").Value = 16861.46
$ws.Range("N74").Value = -18609.46

# Row 77: 'Purified Polyrhythm (L)' / 'Dark Chestnut Rod'
$ws.Range("H77").Value = 16861.46
$ws.Range("J77").Value = 16861.46
$ws.Range("L77").Value = 50584.38
$ws.Range("N77").Value = -59320.38

# Row 107: 'Built to Last' / 'White Oak Lumber'
$ws.Range("H107").Value = 927.2
$ws.Range("I107").Value = 537.1539
$ws.Range("K107").Value = 537.1539
$ws.Range("M107").Value = 1382.8461

# Row 141: 'No Greater Treasure' / 'Claro Walnut Necklace of Gathering'
$ws.Range("H141").Value = 24026.316
$ws.Range("J141").Value = 24026.316
$ws.Range("L141").Value = 24026.316
$ws.Range("N141").Value = -34386.316

$ws = $wb.Worksheets.Item("CUL")
# Row 59: 'Comfort Me with Mushrooms' / 'Buttons in a Blanket'
$ws.Range("H59").Value = 3116
$ws.Range("J59").Value = 3116
$ws.Range("L59").Value = 9348
$ws.Range("N59").Value = -10428

# Row 87: 'Soup That Eats Like a Knight' / 'Clam Chowder'
$ws.Range("H87").Value = 12671.429
$ws.Range("J87").Value = 15950
$ws.Range("L87").Value = 47850
$ws.Range("N87").Value = -50346

# Row 90: 'Like Ma Used to Make (L)' / 'Clam Chowder'
$ws.Range("H90").Value = 12671.429
$ws.Range("J90").Value = 15950
$ws.Range("L90").Value = 143550
$ws.Range("N90").Value = -156030

# Row 92: 'Oh No Udon' / 'Gyr Abanian Flour'
$ws.Range("H92").Value = 1093.4584
$ws.Range("I92").Value = 1434
$ws.Range("J92").Value = 1003.8421
$ws.Range("K92").Value = 4302
$ws.Range("L92").Value = 3011.5263
$ws.Range("M92").Value = -3054
$ws.Range("N92").Value = -5507.5263

# Row 113: 'Can''t Eat Just One' / 'Night Vinegar'
$ws.Range("H113").Value = 1489.5128
$ws.Range("I113").Value = 2443.2666
$ws.Range("J113").Value = 893.4167
$ws.Range("K113").Value = 7329.7998
$ws.Range("L113").Value = 2680.2501
$ws.Range("M113").Value = -5159.7998
$ws.Range("N113").Value = -7020.2501

# Row 118: 'Teetotally' / 'Masala Chai'
$ws.Range("H118").Value = 1753.1666
$ws.Range("I118").Value = 643
$ws.Range("J118").Value = 2863.3333
$ws.Range("K118").Value = 1929
$ws.Range("L118").Value = 8589.999899999999
$ws.Range("M118").Value = -686
$ws.Range("N118").Value = -11075.9999

# Row 122: 'Salt of the North' / 'Northern Sea Salt'
$ws.Range("H122").Value = 1254.6666
$ws.Range("I122").Value = 631.3333
$ws.Range("J122").Value = 1432.762
$ws.Range("K122").Value = 5681.9997
$ws.Range("L122").Value = 12894.858
$ws.Range("M122").Value = -3231.9997
$ws.Range("N122").Value = -17794.858

$ws = $wb.Worksheets.Item("GSM")
# Row 57: 'Gold Is So Last Year' / 'Electrum Circlet (Amber)'
$ws.Range("H57").Value = 13507.625
$ws.Range("J57").Value = 14015.25
$ws.Range("L57").Value = 14015.25
$ws.Range("N57").Value = -15655.25

$ws = $wb.Worksheets.Item("LTW")
# Row 68: 'You Could Say It''s a Moving Target' / 'Wyvern Leather'
$ws.Range("H68").Value = 2334.2354
$ws.Range("I68").Value = 1063
$ws.Range("J68").Value = 8266.666999999999
$ws.Range("K68").Value = 1063
$ws.Range("L68").Value = 8266.666999999999
$ws.Range("M68").Value = -314
$ws.Range("N68").Value = -9764.666999999999

# Row 71: 'They Call It Bloody Mary (L)' / 'Wyvern Leather'
$ws.Range("H71").Value = 2334.2354
$ws.Range("I71").Value = 1063
$ws.Range("J71").Value = 8266.666999999999
$ws.Range("K71").Value = 5315
$ws.Range("L71").Value = 41333.335
$ws.Range("M71").Value = -1571
$ws.Range("N71").Value = -48821.335

# Row 132: 'Tenets of Tanning' / 'Silver Lobo Leather'
$ws.Range("H132").Value = 2458.6667
$ws.Range("I132").Value = 1763.9584
$ws.Range("J132").Value = 3252.6191
$ws.Range("K132").Value = 5291.8752
$ws.Range("L132").Value = 9757.8573
$ws.Range("M132").Value = -2761.8752
$ws.Range("N132").Value = -14817.8573

$ws = $wb.Worksheets.Item("WVR")
# Row 113: 'A Tender Table' / 'Pixie Floss'
$ws.Range("H113").Value = 1411.4166
$ws.Range("I113").Value = 204.85715
$ws.Range("K113").Value = 614.5714499999999
$ws.Range("M113").Value = 1555.42855
